## update scripts with new tpm
# This re-computes the NATMI Vtn -> Tnfrsf11b ligand-receptor table for the
# OldD7 TPM dataset. The "Resolving-Mac" sending-cluster rows are dropped
# (no longer produced by the updated script) while "Resolving-Mac" now
# appears as a target cluster instead, and all numeric columns (expression
# values, specificities, weights) are refreshed with the new TPM-derived
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 11-13 (Resolving-Mac as sending cluster) - they are
# not present anymore in the updated output.
$ws.Rows("11:13").Delete()

# Update rows 2-10 with the new TPM-derived values
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vtn"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 7.844453333333334
$ws.Range("H2").Value = 23.53336
$ws.Range("I2").Value = 0.1489290605659587
$ws.Range("J2").Value = 0.1489290605659588
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 2.214957333333333
$ws.Range("N2").Value = 6.644871999999999
$ws.Range("O2").Value = 0.8812411509483107
$ws.Range("P2").Value = 0.8812411509483107
$ws.Range("Q2").Value = 17.37512943665778
$ws.Range("R2").Value = 156.37616492992
$ws.Range("S2").Value = 0.1312424167427961
$ws.Range("T2").Value = 0.1312424167427962

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vtn"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 7.844453333333334
$ws.Range("H3").Value = 23.53336
$ws.Range("I3").Value = 0.1489290605659587
$ws.Range("J3").Value = 0.1489290605659588
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 0.274148
$ws.Range("N3").Value = 0.8224440000000001
$ws.Range("O3").Value = 0.1090723037479928
$ws.Range("P3").Value = 0.1090723037479928
$ws.Range("Q3").Value = 2.150541192426667
$ws.Range("R3").Value = 19.35487073184
$ws.Range("S3").Value = 0.01624403573095347
$ws.Range("T3").Value = 0.01624403573095347

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vtn"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 7.844453333333334
$ws.Range("H4").Value = 23.53336
$ws.Range("I4").Value = 0.1489290605659587
$ws.Range("J4").Value = 0.1489290605659588
$ws.Range("K4").Value = 1.0
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02434666666666667
$ws.Range("N4").Value = 0.07304
$ws.Range("O4").Value = 0.009686545303696538
$ws.Range("P4").Value = 0.009686545303696536
$ws.Range("Q4").Value = 0.1909862904888889
$ws.Range("R4").Value = 1.7188766144
$ws.Range("S4").Value = 0.001442608092209125
$ws.Range("T4").Value = 0.001442608092209125

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vtn"
$ws.Range("C5").Value = "Tnfrsf11b"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 20.35396833333334
$ws.Range("H5").Value = 61.06190500000001
$ws.Range("I5").Value = 0.3864255740794268
$ws.Range("J5").Value = 0.3864255740794268
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 2.214957333333333
$ws.Range("N5").Value = 6.644871999999999
$ws.Range("O5").Value = 0.8812411509483107
$ws.Range("P5").Value = 0.8812411509483107
$ws.Range("Q5").Value = 45.08317142235111
$ws.Range("R5").Value = 405.7485428011601
$ws.Range("S5").Value = 0.3405341176576158
$ws.Range("T5").Value = 0.3405341176576158

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vtn"
$ws.Range("C6").Value = "Tnfrsf11b"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 20.35396833333334
$ws.Range("H6").Value = 61.06190500000001
$ws.Range("I6").Value = 0.3864255740794268
$ws.Range("J6").Value = 0.3864255740794268
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.274148
$ws.Range("N6").Value = 0.8224440000000001
$ws.Range("O6").Value = 0.1090723037479928
$ws.Range("P6").Value = 0.1090723037479928
$ws.Range("Q6").Value = 5.579999710646668
$ws.Range("R6").Value = 50.21999739582001
$ws.Range("S6").Value = 0.04214832759198374
$ws.Range("T6").Value = 0.04214832759198374

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vtn"
$ws.Range("C7").Value = "Tnfrsf11b"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 20.35396833333334
$ws.Range("H7").Value = 61.06190500000001
$ws.Range("I7").Value = 0.3864255740794268
$ws.Range("J7").Value = 0.3864255740794268
$ws.Range("K7").Value = 1.0
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02434666666666667
$ws.Range("N7").Value = 0.07304
$ws.Range("O7").Value = 0.009686545303696538
$ws.Range("P7").Value = 0.009686545303696536
$ws.Range("Q7").Value = 0.4955512823555556
$ws.Range("R7").Value = 4.4599615412
$ws.Range("S7").Value = 0.00374312882982731
$ws.Range("T7").Value = 0.003743128829827309

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Vtn"
$ws.Range("C8").Value = "Tnfrsf11b"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 24.47399366666667
$ws.Range("H8").Value = 73.421981
$ws.Range("I8").Value = 0.4646453653546145
$ws.Range("J8").Value = 0.4646453653546145
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 2.214957333333333
$ws.Range("N8").Value = 6.644871999999999
$ws.Range("O8").Value = 0.8812411509483107
$ws.Range("P8").Value = 0.8812411509483107
$ws.Range("Q8").Value = 54.20885174793688
$ws.Range("R8").Value = 487.879665731432
$ws.Range("S8").Value = 0.4094646165478988
$ws.Range("T8").Value = 0.4094646165478988

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Vtn"
$ws.Range("C9").Value = "Tnfrsf11b"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 24.47399366666667
$ws.Range("H9").Value = 73.421981
$ws.Range("I9").Value = 0.4646453653546145
$ws.Range("J9").Value = 0.4646453653546145
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.274148
$ws.Range("N9").Value = 0.8224440000000001
$ws.Range("O9").Value = 0.1090723037479928
$ws.Range("P9").Value = 0.1090723037479928
$ws.Range("Q9").Value = 6.709496415729334
$ws.Range("R9").Value = 60.38546774156401
$ws.Range("S9").Value = 0.05067994042505561
$ws.Range("T9").Value = 0.05067994042505562

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Vtn"
$ws.Range("C10").Value = "Tnfrsf11b"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 24.47399366666667
$ws.Range("H10").Value = 73.421981
$ws.Range("I10").Value = 0.4646453653546145
$ws.Range("J10").Value = 0.4646453653546145
$ws.Range("K10").Value = 1.0
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02434666666666667
$ws.Range("N10").Value = 0.07304
$ws.Range("O10").Value = 0.009686545303696538
$ws.Range("P10").Value = 0.009686545303696536
$ws.Range("Q10").Value = 0.5958601658044445
$ws.Range("R10").Value = 5.36274149224
$ws.Range("S10").Value = 0.004500808381660103
$ws.Range("T10").Value = 0.004500808381660103
